# Purchase entry UI changes
#
# 1. Leave the "Gatepass" sheet with a plain B4 selection (no longer the
#    active tab).
# 2. Switch focus to the "PurchaseEntry" sheet, update the medicine/name
#    entries in A2:A4 and B4, and leave the selection on B6.

$wb = $excel.ActiveWorkbook

# Touch "Gatepass" first: move the selection to B4 and drop it as the
# active tab (it was the active sheet before this edit).
$gatepass = $wb.Worksheets.Item("Gatepass")
$gatepass.Activate()
$gatepass.Range("B4").Select()

# Now edit "PurchaseEntry" - this becomes the new active sheet.
$purchaseEntry = $wb.Worksheets.Item("PurchaseEntry")
$purchaseEntry.Activate()

$purchaseEntry.Range("B4").Value = "GOLD COAT LOTION 120ML"
$purchaseEntry.Range("A4").Value = "ne31"
$purchaseEntry.Range("A2").Value = "ne12"
$purchaseEntry.Range("A3").Value = "ne41"

$purchaseEntry.Range("B6").Select()
